$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quantif")

# The "Ladder" volume columns (B, D, E, J, K, L, M) on the Plate1 ladder row (76)
# had been accidentally filled in with values copied from elsewhere; clear them
# back out so the dependent "Volume of ladder" rows below (77:83) correctly
# compute to 0 for those columns again.
$ws.Range("B76").Value = ""
$ws.Range("D76").Value = ""
$ws.Range("E76").Value = ""
$ws.Range("J76").Value = ""
$ws.Range("K76").Value = ""
$ws.Range("L76").Value = ""
$ws.Range("M76").Value = ""

# Restore the window's scroll position / active selection as it was left in
# the saved workbook.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 64
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I71").Select()
